# Update crypto price/volume data per the Tue Apr 18 21:27:40 UTC 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.419.42'
$ws.Range('E2').Value = '  +2.37%  '
$ws.Range('D3').Value = '2.095.13'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('E4').Value = '  -0.94%  '
$ws.Range('D5').Value = '''343.20'
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('D6').Value = '''1.002'
$ws.Range('E6').Value = '  -0.82%  '
$ws.Range('D7').Value = '''0.5236'
$ws.Range('E7').Value = '  +1.54%  '
$ws.Range('D8').Value = '''0.4423'
$ws.Range('E8').Value = '  +1.30%  '
$ws.Range('E9').Value = '  +3.16%  '
$ws.Range('D10').Value = '''0.09330'
$ws.Range('E10').Value = '  +1.17%  '
$ws.Range('D11').Value = '''1.169'
$ws.Range('E11').Value = '  +0.66%  '
$ws.Range('D12').Value = '''24.85'
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('D13').Value = '''8.594'
$ws.Range('E13').Value = '  +4.01%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '''6.903'
$ws.Range('E14').Value = '  +2.48%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.032.75'
$ws.Range('E15').Value = '  -3.12%  '
$ws.Range('D16').Value = '''101.41'
$ws.Range('E16').Value = '  +2.29%  '
$ws.Range('D17').Value = '''0.00001159'
$ws.Range('E17').Value = '  +1.02%  '
$ws.Range('D18').Value = '''1.002'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('D19').Value = '''21.15'
$ws.Range('E19').Value = '  +2.05%  '
$ws.Range('D20').Value = '''0.06663'
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('D21').Value = '''6.328'
$ws.Range('E21').Value = '  +2.43%  '
$ws.Range('D22').Value = '''1.001'
$ws.Range('E22').Value = '  -0.78%  '
$ws.Range('D23').Value = '30.427.12'
$ws.Range('E23').Value = '  +2.34%  '
$ws.Range('D24').Value = '''12.54'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('E25').Value = '  -0.41%  '
$ws.Range('D26').Value = '''21.83'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').Value = '''163.06'
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('D28').Value = '''2.505'
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('D29').Value = '''133.16'
$ws.Range('E29').Value = '  +0.31%  '
$ws.Range('D30').Value = '''1.139'
$ws.Range('E30').Value = '  +1.06%  '
$ws.Range('D31').Value = '''1.662'
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('D32').Value = '''0.1046'
$ws.Range('E32').Value = '  -0.39%  '
$ws.Range('D33').Value = '''6.840'
$ws.Range('E33').Value = '  +9.67%  '
$ws.Range('D34').Value = '''6.261'
$ws.Range('E34').Value = '  +1.92%  '
$ws.Range('D35').Value = '''3.861'
$ws.Range('E35').Value = '  -1.93%  '
$ws.Range('D36').Value = '''10.15'
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('D37').Value = '''0.02640'
$ws.Range('E37').Value = '  +3.28%  '
$ws.Range('D38').Value = '''0.06849'
$ws.Range('E38').Value = '  +2.60%  '
$ws.Range('D39').Value = '''0.6991'
$ws.Range('E39').Value = '  +2.03%  '
$ws.Range('D40').Value = '''12.57'
$ws.Range('E40').Value = '  +1.22%  '
$ws.Range('D41').Value = '''1.344'
$ws.Range('E41').Value = '  +0.19%  '
$ws.Range('E42').Value = '  -0.34%  '
$ws.Range('D43').Value = '''0.6818'
$ws.Range('E43').Value = '  +2.52%  '
$ws.Range('D44').Value = '''14.39'
$ws.Range('E44').Value = '  +1.19%  '
$ws.Range('D45').Value = '''2.345'
$ws.Range('E45').Value = '  +1.61%  '
$ws.Range('D46').Value = '''1.001'
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('E47').Value = '  +18.39%  '
$ws.Range('D48').Value = '''3.638'
$ws.Range('E48').Value = '  +0.43%  '
$ws.Range('D49').Value = '''0.00000000347'
$ws.Range('E49').Value = '  -3.22%  '
$ws.Range('D50').Value = '''1.225'
$ws.Range('E50').Value = '  +10.00%  '
$ws.Range('D51').Value = '''1.214'
$ws.Range('E51').Value = '  -0.20%  '
